# Updated symbol list with latest coinranking.com data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin names / links) -------------------------------
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'

# --- Numeric-looking text cells (price / volume %) ------------------------
# Force these to stay text (matching the source sheet, which stores every
# cell as inline text) by switching to a text number format before the
# write, then resetting the style back to Normal so no formatting residue
# is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '302.90'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.30%'
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.42%'
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.036'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.13%'
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07905'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-2.67%'
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.848'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-5.20%'
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.779'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.20%'
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9177'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.36%'
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1352'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-4.69%'
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1890'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-2.05%'
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09033'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-2.48%'
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03475'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.20%'
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09805'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.45%'
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001401'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.88%'
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006117'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '6.43%'
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.720'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '3.19%'
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.105'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.92%'
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.293'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '10.71%'
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3437'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.05%'
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.186'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '5.90%'
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1309'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.95%'
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-8.96%'
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04409'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-2.13%'
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.54%'
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004603'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-5.51%'
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001302'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '4.92%'
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0004446'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.10%'
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01930'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-4.02%'
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05219'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '5.78%'
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007598'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.67%'
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01018'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-4.31%'
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-3.19%'
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002153'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '2.46%'
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009997'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.00%'
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006154'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-4.40%'
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.12%'
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '64.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '0.45%'
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001661'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '39.50%'
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002102'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.12%'
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002002'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.12%'
$ws.Range("E51").Style = "Normal"
